# This script replicates the existing "ASC1" block (rows 2-30) as a new
# "ASC2" block (rows 32-60), with an extra blank separator row before each
# of the three 9-row sub-blocks (rows 41 and 51), matching the pattern
# already used for the ASC1 block (rows 11 and 21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Constant text (column E/F are identical for every data row in the sheet)
$colE = "MADsFromMedian"
$colF = "latestValue"

# The H/I/J values repeat identically within every 9-row sub-block
$hVals = @(7.5, 3, 1.5, 0.5, -0.5, -1.5, -3, -7.5, -1000)
$iVals = @(1000, 7.5, 3, 1.5, 0.5, -0.5, -1.5, -3, -7.5)
$jVals = @("[Manually written]", "much higher than", "higher than", "slightly higher than", "broadly similar to", "slightly lower than", "lower than", "much lower than", "[Manually written]")

# Per sub-block config: starting row, subject number (column B),
# and the G (column G) text per offset within the block (0-8).
$subjects = @(
    @{ Start = 32; B = 1; G = @("NN_median","NN_median","NN_median","NN_median","NN_median","NN_median","NN_median","NN_median","NN_median") },
    @{ Start = 42; B = 2; G = @("earliestValue","earliestValue","earliestValue","earliestValue","earliestValue","earliestValue","earliestValue","earliestValue","earliestValue") },
    @{ Start = 52; B = 3; G = @("penultimateValue","earliestValue","earliestValue","earliestValue","earliestValue","earliestValue","earliestValue","earliestValue","earliestValue") }
)

foreach ($sub in $subjects) {
    $start = $sub.Start
    $last = $start + 8

    # Copy cell formatting (styles only) from the matching ASC1 sub-block
    # (rows 2-10 / 12-20 / 22-30) so the new rows get identical s="1"/s="2"
    # style indices as the corresponding source rows.
    if ($sub.B -eq 1) { $srcStart = 2 }
    elseif ($sub.B -eq 2) { $srcStart = 12 }
    else { $srcStart = 22 }
    $srcLast = $srcStart + 8
    $ws.Range("A" + $srcStart + ":J" + $srcLast).Copy()
    $ws.Range("A" + $start + ":J" + $last).PasteSpecial(-4122)  # xlPasteFormats

    for ($off = 0; $off -lt 9; $off++) {
        $r = $start + $off
        $ws.Range("A" + $r).Value = "ASC"
        $ws.Range("B" + $r).Value = $sub.B
        $ws.Range("C" + $r).Value = "ASC2"
        $ws.Range("E" + $r).Value = $colE
        $ws.Range("F" + $r).Value = $colF
        $ws.Range("G" + $r).Value = $sub.G[$off]
        $ws.Range("H" + $r).Value = $hVals[$off]
        $ws.Range("I" + $r).Value = $iVals[$off]
        $ws.Range("J" + $r).Value = $jVals[$off]
    }

    # Column D: first row of the block gets its own formula, the remaining
    # 8 rows are filled as one range so they form a shared formula group
    # (matching the existing ASC1 blocks' D2/D3:D10, D12/D13:D20, D22/D23:D30).
    $ws.Range("D" + $start).Formula = "=_xlfn.CONCAT(C" + $start + ",""_"",B" + $start + ")"
    $ws.Range("D" + ($start + 1) + ":D" + $last).Formula = "=_xlfn.CONCAT(C" + ($start + 1) + ",""_"",B" + ($start + 1) + ")"
}

# Blank separator rows (mirroring rows 11 and 21 that separate the ASC1 sub-blocks)
$ws.Range("A11:J11").Copy()
$ws.Range("A41:J41").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A21:J21").Copy()
$ws.Range("A51:J51").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0
$excel.Calculate()

# Update the view: scroll so row 10 is at the top and select M34
# (cosmetic sheetView state captured in the target workbook).
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M34").Select()

Write-Host "ASC2 block added (rows 32-60)."
